$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 37 (Сб. Фитонефрол (Урологический сбор) 50г) - discontinued product
$ws.Rows(37).Delete()

# Step 2: refresh product names and stock quantities for all remaining rows (2-109)
# following the updated stock report; rows were re-sorted by the source system.
$data = @(
  @(2, 'ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы', 63520),
  @(3, 'Солодка корни 50г', 11167),
  @(4, 'Ромашка цветки вн 50г', 59010),
  @(5, 'Мята перечная листья 50г', 14028),
  @(6, 'Череда трава 50г', 12826),
  @(7, 'Шалфей листья 50г', 19388),
  @(8, 'Шиповник плоды низковитаминные 50г', 15100),
  @(9, 'Ламинарии слоевища (морская капуста) 100г', 12265),
  @(10, 'Толокнянка листья 50г', 6282),
  @(11, 'Зверобой трава 50г', 30895),
  @(12, 'Дуба кора 75г', 63360),
  @(13, 'Укроп пахучий плоды 50г', 84925),
  @(14, 'Пижма цветки 75г', 17948),
  @(15, 'Спорыш трава 50г', 14228),
  @(16, 'Крапива листья 50г', 15895),
  @(17, 'Аир корневища 75г', 7300),
  @(18, 'Ноготки цветки 50г', 17710),
  @(19, 'Кукуруза столбики с рыльцами 40г', 31055),
  @(20, 'Сенна листья 50г', 28629),
  @(21, 'Можжевельник плоды 50г', 11634),
  @(22, 'Береза почки 50г', 12572),
  @(23, 'Тысячелистник трава 50г', 16939),
  @(24, 'Лен семена 100г', 85497),
  @(25, 'Крушина кора 50г', 12927),
  @(26, 'Полынь горькая трава 50г', 53501),
  @(27, 'Бессмертник песчаный цветки 30г', 30968),
  @(28, 'Эрва шерстистая трава 30г', 19305),
  @(29, 'Чага (березовый гриб) 50г', 36596),
  @(30, 'Боярышник плоды 75г', 26289),
  @(31, 'Чистотел трава 50г', 25515),
  @(32, 'Валериана корневища с корнями 50г', 20750),
  @(33, 'Эвкалипт прутовидный листья 75г', 14154),
  @(34, 'Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г', 8657),
  @(35, 'Брусника листья 50г', 22739),
  @(36, 'Багульник болотный побеги 50г', 15190),
  @(37, 'Алтей корни 75г', 8869),
  @(38, 'Мать-и-мачеха листья 35г', 21171),
  @(39, 'Чабрец трава 50г', 27266),
  @(40, 'Девясил корневища и корни 50г', 18004),
  @(41, 'Липа цветки 35г', 18799),
  @(42, 'Пустырник трава 50г', 19353),
  @(43, 'Сб. Фитопектол №1 (Грудной сбор №1) 35г', 4610),
  @(44, 'Сб. Грудной №4 50г', 32121),
  @(45, 'Подорожник большой листья 50г', 17622),
  @(46, 'Сб. Фитопектол №2 (Грудной сбор №2) 35г', 7570),
  @(47, 'Рябина плоды 50г', 6202),
  @(48, 'Фп "ФармаЦветик® Фиточай для кормящих мам" 20х1,5 г', 1870),
  @(49, 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г', 3440),
  @(50, 'Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г', 2130),
  @(51, 'Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г', 2350),
  @(52, 'Фп Фиточай "Лактафитол" (БАД) 20х1,5 г', 19994),
  @(53, 'Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г', 4440),
  @(54, 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г', 3970),
  @(55, 'Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г', 2780),
  @(56, 'Фп Брусника листья 20х1,5г', 26609),
  @(57, 'Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г', 16937),
  @(58, 'Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г', 5724),
  @(59, 'Фп Шиповник плоды 20х2,0г', 14850),
  @(60, 'Фп Череда трава 20х1,5г', 27153),
  @(61, 'Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г', 97916),
  @(62, 'Фп Сб. Желудочный №3 20x2,0г', 9864),
  @(63, 'Фп Сенна листья 20x1,5г', 51588),
  @(64, 'Фп Шалфей листья 20х1,5г', 50350),
  @(65, 'Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г', 12591),
  @(66, 'Фп Зверобой трава 20x1,5г', 21090),
  @(67, 'Фп Хвощ полевой трава 20х1,5г', 17601),
  @(68, 'Фп Толокнянка листья 20x1,5г', 17657),
  @(69, 'Фп Ромашка цветки 20x1,5г', 565813),
  @(70, 'Фп Мелисса лекарственная трава 20x1,5г', 20538),
  @(71, 'Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г', 3492),
  @(72, 'Фп Сб. Бруснивер 20x2,0г', 156658),
  @(73, 'Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г', 56814),
  @(74, 'Фп Крапива листья 20x1,5г', 41385),
  @(75, 'Фп Чистотел трава 20х1,5г', 19486),
  @(76, 'Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г', 5526),
  @(77, 'Фп Сб. Арфазетин-Э 20x2,0г', 35416),
  @(78, 'Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г', 42768),
  @(79, 'Фп Мята перечная листья 20x1,5г', 41923),
  @(80, 'Фп Боярышник плоды 20х3,0г', 8532),
  @(81, 'Фп Фиалка трехцветная трава 20x1,5г', 2952),
  @(82, 'Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г', 58507),
  @(83, 'Фп Пижма цветки 20х1,5г', 7194),
  @(84, 'Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г', 3562),
  @(85, 'Фп Пустырник трава 20x1,5г', 24983),
  @(86, 'Фп Пастушья сумка трава 20х1,5г', 7535),
  @(87, 'Фп Золототысячник трава 20х1,5г', 3034),
  @(88, 'Фп Сб. Элекасол 20x2,0г', 19818),
  @(89, 'Фп Чабрец трава 20x1,5 г', 39150),
  @(90, 'Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г', 4824),
  @(91, 'Фп Подорожник листья 20x1,5г', 21897),
  @(92, 'Фп Ольха соплодия 20х1,5г', 4370),
  @(93, 'Фп Душица трава 20x1,5г', 22825),
  @(94, 'Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г', 3150),
  @(95, 'Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г', 3240),
  @(96, 'Фп Береза листья 20x1,5г', 3903),
  @(97, 'Фп Аир корневища 20x1,5г', 5818),
  @(98, 'Фп Липа цветки 20x1,5г', 45030),
  @(99, 'Фп Сб. Грудной №4 20x2,0г', 754992),
  @(100, 'Фп Ноготки цветки 20x1,5г', 30269),
  @(101, 'Фп Почечный чай листья 20x1,5г', 82094),
  @(102, 'Фп Дуб кора 20х1,5г', 25235),
  @(103, 'Фп Кровохлебка корневища и корни 20x1,5г', 10040),
  @(104, 'Фп Крушина кора 20x1,5г', 18948),
  @(105, 'Фп Тысячелистник трава 20x1,5г', 27075),
  @(106, 'Фп Валериана корневища с корнями 20x1,5г', 22424),
  @(107, 'Фп Лапчатка корневища 20x2,5г', 7516),
  @(108, 'Фп Бадан корневища 20x1,5г', 5345),
  @(109, 'Фп Девясил корневища и корни 20х1,5г', 17223)
)

foreach ($item in $data) {
  $r = $item[0]
  $name = $item[1]
  $qty = $item[2]
  $ws.Cells.Item($r, 1).Value = $name
  $ws.Cells.Item($r, 2).Value = $qty
}

# Step 3: normalize number format for cells that previously used the no-thousands-separator style
# (original rows 3, 4, 49, 50; rows 49/50 shift up to 48/49 after the row-37 deletion above)
$ws.Range("B3").NumberFormat = "#,##0"
$ws.Range("B4").NumberFormat = "#,##0"
$ws.Range("B48").NumberFormat = "#,##0"
$ws.Range("B49").NumberFormat = "#,##0"

# Step 4: update selection to match the refreshed view
$ws.Range("A11").Select()